$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Highlight the "Clean log of battle" (C31) and "Token-based security" (C39)
# score cells in yellow, like was done when reviewing progress on the project.
$ws.Range("C31").Interior.Color = 65535
$ws.Range("C39").Interior.Color = 65535

# Clear out the scores for "consider specialities in battle-rounds between
# cards" (C33) and "Contains tracked time" (C51) - not done yet.
$ws.Range("C33").ClearContents()
$ws.Range("C51").ClearContents()

# Update the active selection to reflect where the reviewer left off.
$ws.Range("E36").Select()
